$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()

$ws.Range("B2:C4").Value = "'"
$ws.Range("B2:C4").ClearFormats()
